# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 394
$ws1.Range("F6").Value  = 15
$ws1.Range("F8").Value  = 14045
$ws1.Range("F9").Value  = 114
$ws1.Range("F10").Value = 100
$ws1.Range("F11").Value = 5656
$ws1.Range("F13").Value = 56
$ws1.Range("F15").Value = 52
$ws1.Range("F19").Value = 763
$ws1.Range("F21").Value = 48
$ws1.Range("F22").Value = 10434
$ws1.Range("F24").Value = 34
$ws1.Range("F25").Value = 55
$ws1.Range("F26").Value = 3710
$ws1.Range("F27").Value = 238

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 394
$ws4.Range("F7").Value  = 15
$ws4.Range("F9").Value  = 14045
$ws4.Range("F10").Value = 114
$ws4.Range("F11").Value = 100
$ws4.Range("F12").Value = 5656
$ws4.Range("F14").Value = 56
$ws4.Range("F16").Value = 52
$ws4.Range("F20").Value = 763
$ws4.Range("F22").Value = 48
$ws4.Range("F24").Value = 10434
$ws4.Range("F26").Value = 34
$ws4.Range("F27").Value = 55
$ws4.Range("F28").Value = 3710
$ws4.Range("F29").Value = 238
